$wb = $excel.ActiveWorkbook
$url = "https://magento-demo.mageplaza.com/catalogsearch/result/?q="

# --- Add the new "Search" worksheet after the last existing sheet ---
$accountSheet = $wb.Worksheets.Item("Account")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Search"

# --- Column widths (closest achievable match to source widths - the
#     runtime snaps ColumnWidth to a whole-pixel grid internally) ---
$ws.Columns.Item(1).ColumnWidth = 25.857142857142858
$ws.Columns.Item(2).ColumnWidth = 31.285714285714285
$ws.Columns.Item(3).ColumnWidth = 26.714285714285715
$ws.Columns.Item(4).ColumnWidth = 37.714285714285715

# --- Fill in cell values in the same order the strings were originally authored,
#     so the shared-strings table gets built up in the same order. ---
$ws.Range("A1").Value = "TESTCASEID"
$ws.Range("A5").Value = "Search with no search results"
$ws.Range("C1").Value = "URL_SEARCH_RESULTS"

# C2's hyperlink is added now - it also stamps C2's displayed text with $url
$ws.Hyperlinks.Add($ws.Range("C2"), $url)

$ws.Range("B1").Value = "SEARCH_KEYWORDS"
$ws.Range("B2").Value = "Zing Jump Rope"
$ws.Range("A3").Value = "Search with muti products"
$ws.Range("A4").Value = "Search with muti hints"
$ws.Range("C20").Value = "t"
$ws.Range("A2").Value = "Search with a product"
$ws.Range("B3").Value = "sport"
$ws.Range("B4").Value = "shir"

# --- Remaining values that reuse already-created shared strings ---
$ws.Range("B5").Value = 4567
$ws.Range("C3").Value = $url
$ws.Range("C5").Value = $url
$ws.Hyperlinks.Add($ws.Range("C3:C5"), $url, [Type]::Missing, [Type]::Missing, $url)

# --- Formatting: reuse the workbook's existing named styles by copy/paste-special ---
# Style 1: bold header style (same as row-1 headers on other sheets)
$accountSheet.Range("A1").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)

# Style 2: "Normal 2" vertical-center style used in column A data rows
$accountSheet.Range("A2").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)

# Style 3: Hyperlink style
$accountSheet.Range("D2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("D2:D9").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Sheet view / selection ---
$ws.Range("B5").Select()

# --- Page setup to mirror the other sheets ---
$ws.PageSetup.Orientation = 1

# --- Make "Search" the active (selected) sheet/tab ---
$ws.Activate()
